# InGameChatDialogues.xlsx edit
# Commit: "Update 'what's wrong' dialogue from LV 125 to 131"
#
# Rows 61-74 (column A) hold the dialogue-block's associated level number.
# They were bumped from 125 to 131. Row 61 already carried the "header"
# style (s=9, centered / no explicit font color); rows 62-74 previously
# used the plain centered style (s=1) and now pick up the same s=9 look,
# so copy A61's formatting down over A62:A74 before rewriting the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Match the formatting of A62:A74 to A61 (s="1" -> s="9").
$ws.Range("A61").Copy()
$ws.Range("A62:A74").PasteSpecial(-4122)  # xlPasteFormats

# 2) Bump the level number referenced by this dialogue block: 125 -> 131.
$ws.Range("A61:A74").Value = 131

# 3) Leave the selection where the author's cursor ended up.
[void]$ws.Range("B71").Select()
